# KW-RANK B08CS8YFK5-101-200.xlsx 2020-11-01 22:31:15
# Replace keyword list in column A (rows 1-100) with the new ranked keywords.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'capri hex'
$ws.Range("A2").Value = 'knee pads softball'
$ws.Range("A3").Value = 'wrestling knee pads mens'
$ws.Range("A4").Value = 'snowboarding knee brace'
$ws.Range("A5").Value = 'basketball leggings for boys youth'
$ws.Range("A6").Value = 'honeycomb leggings'
$ws.Range("A7").Value = 'sliding knee pad'
$ws.Range("A8").Value = 'compression knee pad pants'
$ws.Range("A9").Value = 'youth basketball compression pants'
$ws.Range("A10").Value = 'basketball guide hand'
$ws.Range("A11").Value = '28 inch basketball'
$ws.Range("A12").Value = 'mens tights with knee pads'
$ws.Range("A13").Value = 'sliding sleeve softball'
$ws.Range("A14").Value = 'knee sleeve youth wrestling'
$ws.Range("A15").Value = 'basketball pads knee'
$ws.Range("A16").Value = 'sliding pants baseball youth'
$ws.Range("A17").Value = 'sports knee pads basketball'
$ws.Range("A18").Value = 'knee pads crossfit'
$ws.Range("A19").Value = 'boys knee pad leggings'
$ws.Range("A20").Value = 'knee pad basketball'
$ws.Range("A21").Value = 'basketball pads youth'
$ws.Range("A22").Value = 'knee pads for snowboarding'
$ws.Range("A23").Value = 'basketball leg sleeves for boys'
$ws.Range("A24").Value = 'basketball leggings youth boys'
$ws.Range("A25").Value = '20 30 compression leggings'
$ws.Range("A26").Value = 'workout knee pads'
$ws.Range("A27").Value = 'honeycomb compression pants'
$ws.Range("A28").Value = 'boys basketball leg sleeve'
$ws.Range("A29").Value = 'knee brace hockey'
$ws.Range("A30").Value = 'youth thigh compression sleeve'
$ws.Range("A31").Value = 'basketball tights youth boys'
$ws.Range("A32").Value = 'boys volleyball knee pads'
$ws.Range("A33").Value = 'crossfit pads'
$ws.Range("A34").Value = 'compression knee pads youth'
$ws.Range("A35").Value = 'padded basketball leggings for boys'
$ws.Range("A36").Value = 'crossfit knee sleeves men'
$ws.Range("A37").Value = 'baseball pants mens knee high'
$ws.Range("A38").Value = 'basketball knee pads for youth'
$ws.Range("A39").Value = 'softball knee pad'
$ws.Range("A40").Value = 'kneepads men'
$ws.Range("A41").Value = 'knee compression sleeve reduce strain'
$ws.Range("A42").Value = 'knee brace'
$ws.Range("A43").Value = 'knee compression sleeve'
$ws.Range("A44").Value = 'knee support'
$ws.Range("A45").Value = 'knee brace support'
$ws.Range("A46").Value = 'knee brace compression sleeve'
$ws.Range("A47").Value = 'compression sleeve'
$ws.Range("A48").Value = 'knee sleeve'
$ws.Range("A49").Value = 'knee compression'
$ws.Range("A50").Value = 'compression knee brace'
$ws.Range("A51").Value = 'knee pad'
$ws.Range("A52").Value = 'thigh compression sleeve'
$ws.Range("A53").Value = 'basketball training'
$ws.Range("A54").Value = 'crossfit training'
$ws.Range("A55").Value = 'basketball thigh pads'
$ws.Range("A56").Value = 'waist training leggings'
$ws.Range("A57").Value = 'work out tights mens'
$ws.Range("A58").Value = 'compression leggings capri'
$ws.Range("A59").Value = 'black workout leggings'
$ws.Range("A60").Value = 'workout legging'
$ws.Range("A61").Value = 'capri legging'
$ws.Range("A62").Value = 'workout tights men'
$ws.Range("A63").Value = 'workout pad'
$ws.Range("A64").Value = 'leggings basketball'
$ws.Range("A65").Value = 'black legging'
$ws.Range("A66").Value = 'basketball knee sleeve men'
$ws.Range("A67").Value = 'basketball pads for training'
$ws.Range("A68").Value = 'workout tights for men'
$ws.Range("A69").Value = 'knee pads basketball mcdavid black'
$ws.Range("A70").Value = 'basketball knee pads adult black'
$ws.Range("A71").Value = 'mens leggings compression'
$ws.Range("A72").Value = 'black capri legging'
$ws.Range("A73").Value = 'man leggings'
$ws.Range("A74").Value = 'basketball sleeve black'
$ws.Range("A75").Value = 'men workout pants'
$ws.Range("A76").Value = 'basketball knee pad'
$ws.Range("A77").Value = 'hip workout'
$ws.Range("A78").Value = 'compression legging'
$ws.Range("A79").Value = 'workout pant'
$ws.Range("A80").Value = 'capri pant'
$ws.Range("A81").Value = 'youth training basketball'
$ws.Range("A82").Value = 'compression pant'
$ws.Range("A83").Value = 'baseball pants knee high mens'
$ws.Range("A84").Value = 'basketball knee pads youth girls'
$ws.Range("A85").Value = 'basketball knee protector'
$ws.Range("A86").Value = 'basketball leggings with knee pads'
$ws.Range("A87").Value = 'basketball youth knee pads'
$ws.Range("A88").Value = 'compression knee pants'
$ws.Range("A89").Value = 'compression leggings boys basketball'
$ws.Range("A90").Value = 'knee guards for basketball'
$ws.Range("A91").Value = 'knee pad basketball men'
$ws.Range("A92").Value = 'knee pad pants basketball'
$ws.Range("A93").Value = 'knee pads basketball'
$ws.Range("A94").Value = 'knee pads boys basketball'
$ws.Range("A95").Value = 'knee pads compression pants'
$ws.Range("A96").Value = 'knee pads for basketball men'
$ws.Range("A97").Value = 'knee pads men basketball'
$ws.Range("A98").Value = 'knee pads pants for men'
$ws.Range("A99").Value = 'knee pads wrestling'
$ws.Range("A100").Value = 'knee protection for workout'
